# Updated via Streamlit Approval System
# Appends a new pending-approval row (row 12) to the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 12

$ws.Cells.Item($row, 1).Value  = "WGG 02"
$ws.Cells.Item($row, 2).Value  = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item($row, 3).Value  = "15-01-2026"
$ws.Cells.Item($row, 4).Value  = 286962
$ws.Cells.Item($row, 5).Value  = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item($row, 6).Value  = 34413429360
$ws.Cells.Item($row, 7).Value  = "NEFT"
$ws.Cells.Item($row, 8).Value  = "SBIN0003229"
$ws.Cells.Item($row, 9).Value  = "AAAFW8862C"
$ws.Cells.Item($row, 10).Value = "32AAAFW8862C1Z9"
# Column K (BENEFICIARY NAME) intentionally left blank for this row.
$ws.Cells.Item($row, 12).Value = "f2496d5d-1d54-4799-a952-fcc3cbbc08cc"
# Columns M-T left blank for this row.
$ws.Cells.Item($row, 21).Value = "pending"
$ws.Cells.Item($row, 22).Value = 126000
# Column W (FINAL AMOUNT) left blank for this row.
$ws.Cells.Item($row, 24).Value = "Kolkata RPA_UNIQUE_ID : 52342bcc-c106-4bb9-9695-3bf75c83ca21"
$ws.Cells.Item($row, 25).Value = "Kolkata"
$ws.Cells.Item($row, 26).Value = "PAYMENT"
$ws.Cells.Item($row, 27).Value = "accountsassist@westernidc.com"
$ws.Cells.Item($row, 28).Value = "ESTIMATION NOT MATCHED"
$ws.Cells.Item($row, 29).Value = 0
$ws.Cells.Item($row, 30).Value = 0
$ws.Cells.Item($row, 31).Value = 0
# Columns AF-AO left blank for this row.
